# ===========================================================================
# feat: add 2022-Q1 data
#
# The workbook tracks, per fiscal quarter, which mutual funds hold the
# stock and a "总计" (totals) roll-up sheet. This adds a new "2022-Q1"
# sheet (inserted right before "总计") with that quarter's per-fund detail,
# and prepends a matching summary row to "总计".
# ===========================================================================

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Helper: write a string into a cell, preserving it as TEXT. Excel (like
# the real COM object model) auto-coerces a bare numeric-looking string
# ("010336", "44.28", ...) to a number on assignment, which would destroy
# leading zeros / intentionally-textual numeric fields. A leading
# apostrophe forces literal text entry, exactly like a user typing
# '010336 into a cell — but we only do that when the text actually looks
# numeric, so plain labels don't pick up a spurious quote-prefix style.
# ---------------------------------------------------------------------
function Set-Text {
    param($cell, [string]$text)
    if ($text -match '^-?\d+(\.\d+)?$') {
        $cell.Value = "'" + $text
    } else {
        $cell.Value = $text
    }
}

# ---------------------------------------------------------------------
# Step 1: Restructure the sheets.
#
# Today "总计" is the last sheet (sheetId 6). We rename it in place to
# "2022-Q1" (so it keeps sheetId 6 / its existing worksheet part) and
# then add a brand new sheet right after it, named "总计" (this new sheet
# naturally gets the next sheetId, 7) — giving the same sheetId/rId
# assignment the diff shows: 2022-Q1 keeps id 6, the (new) 总计 gets id 7.
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Name = "2022-Q1"

$newTotalSheet = $wb.Worksheets.Add($null, $totalSheet)
$newTotalSheet.Name = "总计"

$ws = $totalSheet
$ws2 = $newTotalSheet

# ---------------------------------------------------------------------
# Step 2: Grab correctly-formatted source cells to copy styling from
# (bold + centered + top-aligned + thin-bordered — the "s=2" style used
# for every sheet's header row and row-index column) before we start
# overwriting "2022-Q1" (née 总计)'s old contents.
# ---------------------------------------------------------------------
$styleSrc = $wb.Worksheets.Item("2021-Q4")
$headerStyleRange = $styleSrc.Range("B1:H1")
$indexStyleCell = $styleSrc.Range("A2")

# ---------------------------------------------------------------------
# Step 3: Clear out the old totals-table content from what is now the
# "2022-Q1" sheet, then lay down the per-fund holdings table.
# ---------------------------------------------------------------------
$ws.Cells.Clear()

# Header row 1 (B1:H1), bold/centered style copied from another sheet.
$headerStyleRange.Copy()
$ws.Range("B1:H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Header labels
# ---- Header row (row 1) ----
Set-Text $ws.Cells.Item(1, 2) '基金代码'
Set-Text $ws.Cells.Item(1, 3) '基金名称'
Set-Text $ws.Cells.Item(1, 4) '基金规模'
Set-Text $ws.Cells.Item(1, 5) '股票总仓位'
Set-Text $ws.Cells.Item(1, 6) '仓位占比'
Set-Text $ws.Cells.Item(1, 7) '持有市值(亿元)'
Set-Text $ws.Cells.Item(1, 8) '仓位排名'

# Row-index column (A2:A19), bold/centered style copied from another sheet.
$indexStyleCell.Copy()
$ws.Range("A2:A19").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data rows 2..19
$ws.Cells.Item(2, 1).Value = 0
Set-Text $ws.Cells.Item(2, 2) '010336'
Set-Text $ws.Cells.Item(2, 3) '中欧悦享生活混合A'
Set-Text $ws.Cells.Item(2, 4) '44.28'
Set-Text $ws.Cells.Item(2, 5) '90.44'
Set-Text $ws.Cells.Item(2, 6) '7.86'
Set-Text $ws.Cells.Item(2, 7) '3.4804'
$ws.Cells.Item(2, 8).Value = 5
$ws.Cells.Item(3, 1).Value = 1
Set-Text $ws.Cells.Item(3, 2) '010680'
Set-Text $ws.Cells.Item(3, 3) '华夏新兴成长股票A'
Set-Text $ws.Cells.Item(3, 4) '55.80'
Set-Text $ws.Cells.Item(3, 5) '87.37'
Set-Text $ws.Cells.Item(3, 6) '5.76'
Set-Text $ws.Cells.Item(3, 7) '3.2141'
$ws.Cells.Item(3, 8).Value = 4
$ws.Cells.Item(4, 1).Value = 2
Set-Text $ws.Cells.Item(4, 2) '010305'
Set-Text $ws.Cells.Item(4, 3) '华夏创新驱动混合A'
Set-Text $ws.Cells.Item(4, 4) '28.56'
Set-Text $ws.Cells.Item(4, 5) '88.57'
Set-Text $ws.Cells.Item(4, 6) '5.73'
Set-Text $ws.Cells.Item(4, 7) '1.6365'
$ws.Cells.Item(4, 8).Value = 4
$ws.Cells.Item(5, 1).Value = 3
Set-Text $ws.Cells.Item(5, 2) '002621'
Set-Text $ws.Cells.Item(5, 3) '中欧消费主题股票A'
Set-Text $ws.Cells.Item(5, 4) '19.29'
Set-Text $ws.Cells.Item(5, 5) '88.29'
Set-Text $ws.Cells.Item(5, 6) '8.29'
Set-Text $ws.Cells.Item(5, 7) '1.5991'
$ws.Cells.Item(5, 8).Value = 4
$ws.Cells.Item(6, 1).Value = 4
Set-Text $ws.Cells.Item(6, 2) '002229'
Set-Text $ws.Cells.Item(6, 3) '华夏经济转型股票'
Set-Text $ws.Cells.Item(6, 4) '11.71'
Set-Text $ws.Cells.Item(6, 5) '86.14'
Set-Text $ws.Cells.Item(6, 6) '5.47'
Set-Text $ws.Cells.Item(6, 7) '0.6405'
$ws.Cells.Item(6, 8).Value = 2
$ws.Cells.Item(7, 1).Value = 5
Set-Text $ws.Cells.Item(7, 2) '002697'
Set-Text $ws.Cells.Item(7, 3) '中欧消费主题股票C'
Set-Text $ws.Cells.Item(7, 4) '6.20'
Set-Text $ws.Cells.Item(7, 5) '88.29'
Set-Text $ws.Cells.Item(7, 6) '8.29'
Set-Text $ws.Cells.Item(7, 7) '0.5140'
$ws.Cells.Item(7, 8).Value = 4
$ws.Cells.Item(8, 1).Value = 6
Set-Text $ws.Cells.Item(8, 2) '010852'
Set-Text $ws.Cells.Item(8, 3) '中欧内需成长混合型证券投资基金A'
Set-Text $ws.Cells.Item(8, 4) '5.23'
Set-Text $ws.Cells.Item(8, 5) '91.46'
Set-Text $ws.Cells.Item(8, 6) '6.61'
Set-Text $ws.Cells.Item(8, 7) '0.3457'
$ws.Cells.Item(8, 8).Value = 4
$ws.Cells.Item(9, 1).Value = 7
Set-Text $ws.Cells.Item(9, 2) '005620'
Set-Text $ws.Cells.Item(9, 3) '中欧品质消费股票A'
Set-Text $ws.Cells.Item(9, 4) '3.74'
Set-Text $ws.Cells.Item(9, 5) '90.47'
Set-Text $ws.Cells.Item(9, 6) '8.20'
Set-Text $ws.Cells.Item(9, 7) '0.3067'
$ws.Cells.Item(9, 8).Value = 4
$ws.Cells.Item(10, 1).Value = 8
Set-Text $ws.Cells.Item(10, 2) '010681'
Set-Text $ws.Cells.Item(10, 3) '华夏新兴成长股票C'
Set-Text $ws.Cells.Item(10, 4) '5.29'
Set-Text $ws.Cells.Item(10, 5) '87.37'
Set-Text $ws.Cells.Item(10, 6) '5.76'
Set-Text $ws.Cells.Item(10, 7) '0.3047'
$ws.Cells.Item(10, 8).Value = 4
$ws.Cells.Item(11, 1).Value = 9
Set-Text $ws.Cells.Item(11, 2) '006868'
Set-Text $ws.Cells.Item(11, 3) '华夏科技成长股票'
Set-Text $ws.Cells.Item(11, 4) '5.23'
Set-Text $ws.Cells.Item(11, 5) '87.33'
Set-Text $ws.Cells.Item(11, 6) '5.77'
Set-Text $ws.Cells.Item(11, 7) '0.3018'
$ws.Cells.Item(11, 8).Value = 4
$ws.Cells.Item(12, 1).Value = 10
Set-Text $ws.Cells.Item(12, 2) '004350'
Set-Text $ws.Cells.Item(12, 3) '汇丰晋信价值先锋股票'
Set-Text $ws.Cells.Item(12, 4) '4.99'
Set-Text $ws.Cells.Item(12, 5) '93.32'
Set-Text $ws.Cells.Item(12, 6) '2.83'
Set-Text $ws.Cells.Item(12, 7) '0.1412'
$ws.Cells.Item(12, 8).Value = 8
$ws.Cells.Item(13, 1).Value = 11
Set-Text $ws.Cells.Item(13, 2) '010306'
Set-Text $ws.Cells.Item(13, 3) '华夏创新驱动混合C'
Set-Text $ws.Cells.Item(13, 4) '2.15'
Set-Text $ws.Cells.Item(13, 5) '88.57'
Set-Text $ws.Cells.Item(13, 6) '5.73'
Set-Text $ws.Cells.Item(13, 7) '0.1232'
$ws.Cells.Item(13, 8).Value = 4
$ws.Cells.Item(14, 1).Value = 12
Set-Text $ws.Cells.Item(14, 2) '005621'
Set-Text $ws.Cells.Item(14, 3) '中欧品质消费股票C'
Set-Text $ws.Cells.Item(14, 4) '1.11'
Set-Text $ws.Cells.Item(14, 5) '90.47'
Set-Text $ws.Cells.Item(14, 6) '8.20'
Set-Text $ws.Cells.Item(14, 7) '0.0910'
$ws.Cells.Item(14, 8).Value = 4
$ws.Cells.Item(15, 1).Value = 13
Set-Text $ws.Cells.Item(15, 2) '010337'
Set-Text $ws.Cells.Item(15, 3) '中欧悦享生活混合C'
Set-Text $ws.Cells.Item(15, 4) '1.08'
Set-Text $ws.Cells.Item(15, 5) '90.44'
Set-Text $ws.Cells.Item(15, 6) '7.86'
Set-Text $ws.Cells.Item(15, 7) '0.0849'
$ws.Cells.Item(15, 8).Value = 5
$ws.Cells.Item(16, 1).Value = 14
Set-Text $ws.Cells.Item(16, 2) '010853'
Set-Text $ws.Cells.Item(16, 3) '中欧内需成长混合型证券投资基金C'
Set-Text $ws.Cells.Item(16, 4) '0.67'
Set-Text $ws.Cells.Item(16, 5) '91.46'
Set-Text $ws.Cells.Item(16, 6) '6.61'
Set-Text $ws.Cells.Item(16, 7) '0.0443'
$ws.Cells.Item(16, 8).Value = 4
$ws.Cells.Item(17, 1).Value = 15
Set-Text $ws.Cells.Item(17, 2) '004351'
Set-Text $ws.Cells.Item(17, 3) '汇丰晋信珠三角区域发展混合'
Set-Text $ws.Cells.Item(17, 4) '0.51'
Set-Text $ws.Cells.Item(17, 5) '93.92'
Set-Text $ws.Cells.Item(17, 6) '5.15'
Set-Text $ws.Cells.Item(17, 7) '0.0263'
$ws.Cells.Item(17, 8).Value = 4
$ws.Cells.Item(18, 1).Value = 16
Set-Text $ws.Cells.Item(18, 2) '012461'
Set-Text $ws.Cells.Item(18, 3) '西藏东财国证龙头家电指数型发起式证券投资基金A'
Set-Text $ws.Cells.Item(18, 4) '1.01'
Set-Text $ws.Cells.Item(18, 5) '94.99'
Set-Text $ws.Cells.Item(18, 6) '2.37'
Set-Text $ws.Cells.Item(18, 7) '0.0239'
$ws.Cells.Item(18, 8).Value = 9
$ws.Cells.Item(19, 1).Value = 17
Set-Text $ws.Cells.Item(19, 2) '012462'
Set-Text $ws.Cells.Item(19, 3) '西藏东财国证龙头家电指数型发起式证券投资基金C'
Set-Text $ws.Cells.Item(19, 4) '0.42'
Set-Text $ws.Cells.Item(19, 5) '94.99'
Set-Text $ws.Cells.Item(19, 6) '2.37'
Set-Text $ws.Cells.Item(19, 7) '0.0100'
$ws.Cells.Item(19, 8).Value = 9

# ---------------------------------------------------------------------
# Step 4: Populate the new "总计" sheet with the updated roll-up table —
# same shape as before, with a new 2022-Q1 row inserted at the top and
# every other row shifted down by one.
# ---------------------------------------------------------------------
$headerStyleRange.Copy()
$ws2.Range("B1:D1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---- Header row (row 1) ----
Set-Text $ws2.Cells.Item(1, 2) '日期'
Set-Text $ws2.Cells.Item(1, 3) '持有数量(只)'
Set-Text $ws2.Cells.Item(1, 4) '持有市值(亿元)'

$indexStyleCell.Copy()
$ws2.Range("A2:A7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data rows 2..7
$ws2.Cells.Item(2, 1).Value = 0
Set-Text $ws2.Cells.Item(2, 2) '2022-Q1'
$ws2.Cells.Item(2, 3).Value = 18
$ws2.Cells.Item(2, 4).Value = 12.89
$ws2.Cells.Item(3, 1).Value = 1
Set-Text $ws2.Cells.Item(3, 2) '2021-Q4'
$ws2.Cells.Item(3, 3).Value = 22
$ws2.Cells.Item(3, 4).Value = 17.48
$ws2.Cells.Item(4, 1).Value = 2
Set-Text $ws2.Cells.Item(4, 2) '2021-Q3'
$ws2.Cells.Item(4, 3).Value = 12
$ws2.Cells.Item(4, 4).Value = 9.44
$ws2.Cells.Item(5, 1).Value = 3
Set-Text $ws2.Cells.Item(5, 2) '2021-Q2'
$ws2.Cells.Item(5, 3).Value = 14
$ws2.Cells.Item(5, 4).Value = 17.91
$ws2.Cells.Item(6, 1).Value = 4
Set-Text $ws2.Cells.Item(6, 2) '2021-Q1'
$ws2.Cells.Item(6, 3).Value = 28
$ws2.Cells.Item(6, 4).Value = 24.58
$ws2.Cells.Item(7, 1).Value = 5
Set-Text $ws2.Cells.Item(7, 2) '2020-Q4'
$ws2.Cells.Item(7, 3).Value = 35
$ws2.Cells.Item(7, 4).Value = 21.57

Write-Output "done"
